# "Generate Report for Archive"
#
# The localization-status report is regenerated: every "Ready for handoff"
# status cell becomes "In Translation", and the (now shorter) Status /
# zh-cn / de-de columns are narrowed to fit the new text on the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Flip the status text everywhere it appears -------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Narrow the columns that held the status text ------------------------
# Overview: zh-cn status (E) and de-de status (F)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de detail sheets: Status (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
